$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.486.63'
$ws.Range('E2').Value = '  +1.92%  '
$ws.Range('D3').Value = '''1.919.07'
$ws.Range('E3').Value = '  +1.62%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''326.27'
$ws.Range('E5').Value = '  -1.67%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').Value = '''0.4747'
$ws.Range('E7').Value = '  +2.60%  '
$ws.Range('D8').Value = '''0.4099'
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').Value = '''47.88'
$ws.Range('D10').Value = '''0.08057'
$ws.Range('E10').Value = '  +0.86%  '
$ws.Range('E11').Value = '  +1.97%  '
$ws.Range('D12').Value = '''22.54'
$ws.Range('E12').Value = '  +3.74%  '
$ws.Range('D13').Value = '''1.917.51'
$ws.Range('E13').Value = '  +1.79%  '
$ws.Range('D14').Value = '''5.939'
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('D15').Value = '''7.166'
$ws.Range('E15').Value = '  +1.39%  '
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = '''0.00001034'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('E19').Value = '  +0.57%  '
$ws.Range('D20').Value = '''17.79'
$ws.Range('E20').Value = '  +1.79%  '
$ws.Range('D21').Value = '''1.001'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').Value = '''29.489.60'
$ws.Range('E22').Value = '  +1.80%  '
$ws.Range('E23').Value = '  +3.10%  '
$ws.Range('E24').Value = '  +2.29%  '
$ws.Range('D25').Value = '''2.209'
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('D26').Value = '''2.137.25'
$ws.Range('E26').Value = '  +1.33%  '
$ws.Range('D27').Value = '''154.90'
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('D28').Value = '''19.86'
$ws.Range('E28').Value = '  +0.88%  '
$ws.Range('D29').Value = '''6.011'
$ws.Range('E29').Value = '  +11.03%  '
$ws.Range('D30').Value = '''2.133'
$ws.Range('E30').Value = '  +0.72%  '
$ws.Range('D31').Value = '''117.94'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').Value = '''1.066'
$ws.Range('E32').Value = '  +9.14%  '
$ws.Range('D33').Value = '''0.09552'
$ws.Range('E33').Value = '  +2.21%  '
$ws.Range('D34').Value = '''1.435'
$ws.Range('D35').Value = '''3.561'
$ws.Range('E35').Value = '  -1.18%  '
$ws.Range('D36').Value = '''5.419'
$ws.Range('E36').Value = '  +2.66%  '
$ws.Range('E37').Value = '  +0.96%  '
$ws.Range('D38').Value = '''0.02263'
$ws.Range('E38').Value = '  +1.47%  '
$ws.Range('D39').Value = '''8.335'
$ws.Range('E39').Value = '  +0.64%  '
$ws.Range('D40').Value = '''1.175'
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').Value = '''0.5900'
$ws.Range('E41').Value = '  +2.17%  '
$ws.Range('D42').Value = '''2.551'
$ws.Range('E42').Value = '  +11.96%  '
$ws.Range('D43').Value = '''0.1847'
$ws.Range('E43').Value = '  +1.56%  '
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').Value = '''0.08018'
$ws.Range('E45').Value = '  +14.40%  '
$ws.Range('D46').Value = '''1.287'
$ws.Range('E46').Value = '  +2.01%  '
$ws.Range('D47').Value = '''0.5567'
$ws.Range('E47').Value = '  +1.55%  '
$ws.Range('D48').Value = '''12.10'
$ws.Range('E48').Value = '  +0.91%  '
$ws.Range('D49').Value = '''1.938'
$ws.Range('E49').Value = '  +1.60%  '
$ws.Range('D50').Value = '''113.12'
$ws.Range('D51').Value = '''45.13'
$ws.Range('E51').Value = '  +2.07%  '
